# Commit message: "Fruta / hortaliza, semanal"
#
# A new weekly price record is inserted at row 137 of the sheet (pushing the
# existing data for the "Ajo" / Macroferia Regional de Talca series down by
# one row, through the former last data row 273 which becomes row 274).
#
# The new row reuses the constant columns of the series (Mercado ID, Mercado,
# Region, Codreg, Categoria ID/Categoria/Variedad/Calidad, Kg o Unidades,
# Clasificacion) and carries its own date / volume / price data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 137; this shifts rows 137:273 down to 138:274
# and keeps all earlier rows (and their formatting) untouched.
$ws.Rows("137:137").Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A137").Value = 5
$ws.Range("B137").Value = "Macroferia Regional de Talca"
$ws.Range("C137").Value = "Maule"
$ws.Range("D137").Value = 44658
$ws.Range("E137").Value = 7
$ws.Range("F137").Value = 100112003
$ws.Range("G137").Value = "Ajo"
$ws.Range("H137").Value = "Chino"
$ws.Range("I137").Value = "Primera"
$ws.Range("J137").Value = 220
$ws.Range("K137").Value = 21000
$ws.Range("L137").Value = 21000
$ws.Range("M137").Value = 21000
$ws.Range("N137").Value = "`$/malla 10 kilos"
$ws.Range("O137").Value = "China"
$ws.Range("P137").Value = 2100
$ws.Range("Q137").Value = 10
$ws.Range("R137").Value = "Hortaliza"

# Match the date style (yyyy-mm-dd ...) used by the rest of column D.
$ws.Range("D137").NumberFormat = $ws.Range("D138").NumberFormat
